$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2: Jalen Brunson -> Jrue Holiday
$ws.Range("A2").Value = "Jrue Holiday"
$ws.Range("B2").Value = "PG,SG"
$ws.Range("C2").Value = "Boston Celtics"

# Row 3: Trae Young -> Yves Missi
$ws.Range("A3").Value = "Yves Missi"
$ws.Range("B3").Value = "C"
$ws.Range("C3").Value = "New Orleans Pelicans"

# Row 4: Draymond Green -> Alperen Sengün
$ws.Range("A4").Value = "Alperen Sengün"
$ws.Range("B4").Value = "C"
$ws.Range("C4").Value = "Houston Rockets"

# Row 5: Alperen Sengün -> Draymond Green
$ws.Range("A5").Value = "Draymond Green"
$ws.Range("B5").Value = "PF,C"
$ws.Range("C5").Value = "Golden State Warriors"

# Row 6: Dereck Lively II -> Jalen Brunson
$ws.Range("A6").Value = "Jalen Brunson"
$ws.Range("B6").Value = "PG"
$ws.Range("C6").Value = "New York Knicks"

# Row 7: Walker Kessler -> Trae Young
$ws.Range("A7").Value = "Trae Young"
$ws.Range("B7").Value = "PG"
$ws.Range("C7").Value = "Atlanta Hawks"

# Row 8: Yves Missi -> Walker Kessler
$ws.Range("A8").Value = "Walker Kessler"
$ws.Range("B8").Value = "C"
$ws.Range("C8").Value = "Utah Jazz"

# Row 16: Jrue Holiday -> D'Angelo Russell
$ws.Range("A16").Value = "D'Angelo Russell"
$ws.Range("B16").Value = "PG"
$ws.Range("C16").Value = "Los Angeles Lakers"
